$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 60
$ws.Range("B2").Value = 59
$ws.Range("C2").Value = 1

$ws.Range("B5").Value = 0.9833333333333333
$ws.Range("C5").Value = 0.01666666666666667
